# v1.4: monto numerico en exportacion Excel; regla bitacora por defecto
# Applies the changes described by the commit diff to Bitacora_tareas.xlsx

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1) Sheet "Log": append two new bitacora rows (28, 29)
# ---------------------------------------------------------------------------
$log = $wb.Worksheets.Item("Log")

$log.Range("A28").Value = "27/02/2025"
$log.Range("B28").Value = "18:10"
$log.Range("C28").Value = "Solapa Errores global y exportación a Excel"
$log.Range("D28").Value = "Nueva pestaña Errores en el dashboard (a la derecha de Sin cotización) que lista todos los egresos con error de clasificación, permite editar cada registro con el mismo modal de edición y se puede exportar a Excel con todos los campos relevantes (incluyendo editado y editado_detalle)."
$log.Range("E28").Value = "Diagnostico"

$log.Range("A29").Value = "27/02/2025"
$log.Range("B29").Value = "18:20"
$log.Range("C29").Value = "Monto numérico en exportación Excel"
$log.Range("D29").Value = "En ambas exportaciones (Transacciones y Errores), la columna monto se escribe como valor numérico (Number) en lugar de texto, para que Excel reconozca números y permita usar fórmulas (SUM, SUMIF, etc.)."
$log.Range("E29").Value = "Diagnostico"

# ---------------------------------------------------------------------------
# 2) Sheet "Resumen": update the "Exportar a Excel" description (B25)
# ---------------------------------------------------------------------------
$resumen = $wb.Worksheets.Item("Resumen")

$resumen.Range("B25").Value = "Botón en la barra de la tabla (solo icono). Exporta la tabla de transacciones tal como está en Supabase: una hoja ""Transacciones"" con columnas fecha, mes, anio, tipo_movimiento, monto (valor numérico para fórmulas), status, medio_pago, moneda, descripcion, cliente, categoria, cat_desc, origen_archivo, cuenta_contable, editado, editado_detalle. Export Errores: monto también como número. Permite analizar y usar fórmulas en Excel."

# ---------------------------------------------------------------------------
# 3) Sheet "Versiones": append version 1.4 (row 6)
# ---------------------------------------------------------------------------
$versiones = $wb.Worksheets.Item("Versiones")

$versiones.Range("A6").Value = "1.4"
$versiones.Range("B6").Value = "27/02/2025"
$versiones.Range("C6").Value = "Exportación Excel: monto como valor numérico (fórmulas en Excel); regla bitácora por defecto reforzada"

# ---------------------------------------------------------------------------
# 4) New sheet "Presupuesto" (added at the end of the workbook)
# ---------------------------------------------------------------------------
$sheetCount = $wb.Worksheets.Count
$lastSheet = $wb.Worksheets.Item($sheetCount)
$presupuesto = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)
$presupuesto.Name = "Presupuesto"

$presupuesto.Range("A1").Value = "Grupo"
$presupuesto.Range("B1").Value = "Descripción comercial"
$presupuesto.Range("C1").Value = "Importe sugerido (ARS)"

$presupuesto.Range("A2").Value = "Normalización de datos"
$presupuesto.Range("B2").Value = "Relevamiento, limpieza y normalización de datos históricos de caja (campos de moneda, categorías, cuentas contables, flags de edición). Incluye lógica de excepciones y detección de inconsistencias."
$presupuesto.Range("C2").Value = 250000

$presupuesto.Range("A3").Value = "Dashboard flujo de caja"
$presupuesto.Range("B3").Value = "Diseño y desarrollo del dashboard mensual (Flujo por mes, Resumen, alertas, modal By Categoría / By Cuenta, gráficos de serie mensual). Incluye formatos de moneda y visualizaciones."
$presupuesto.Range("C3").Value = 320000

$presupuesto.Range("A4").Value = "Bitácora y documentación"
$presupuesto.Range("B4").Value = "Implementación de la bitácora en Excel (Log, Resumen, Versiones, Ref Git y Vercel, Presupuesto) y documentación funcional básica para el uso de la app."
$presupuesto.Range("C4").Value = 120000

$presupuesto.Range("A5").Value = "Integración y despliegue"
$presupuesto.Range("B5").Value = "Configuración de repositorio Git/GitHub, flujo de despliegue a Vercel y ajustes de configuración (vercel.json, conexión con Supabase)."
$presupuesto.Range("C5").Value = 90000

$presupuesto.Range("A6").Value = "Mantenimiento y soporte inicial"
$presupuesto.Range("B6").Value = "Soporte post–implementación, pequeños ajustes funcionales y acompañamiento durante el primer período de uso."
$presupuesto.Range("C6").Value = 80000

# Column widths matching the authored sheet (A, B, C)
$presupuesto.Columns.Item(1).ColumnWidth = 32.83203125
$presupuesto.Columns.Item(2).ColumnWidth = 90.83203125
$presupuesto.Columns.Item(3).ColumnWidth = 24.83203125
